# Generate Report for Archive
#
# - Update the "Ready for handoff" status text to "In Translation" everywhere
#   it appears (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2 all share the
#   same underlying text).
# - Narrow the "zh-cn"/"de-de" status columns (Overview cols E & F, and the
#   Status column (col C) on the zh-cn / de-de detail sheets).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Update status text wherever it shows "Ready for handoff" -> "In Translation"
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws2.Range("C2").Value = "In Translation"
$ws3.Range("C2").Value = "In Translation"

# Narrow the status columns
$ws1.Columns.Item(5).ColumnWidth = 12.48
$ws1.Columns.Item(6).ColumnWidth = 12.48
$ws2.Columns.Item(3).ColumnWidth = 12.48
$ws3.Columns.Item(3).ColumnWidth = 12.48
